$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 holds the "riparian" habitat_attribute record.
# B7 previously referenced the shared string "120"; it now needs to
# reference the shared string "0" (kept as text, not a number), matching
# the other EDT_greater_than_RTT / RTT_greater_than_EDT columns that use
# the text value "0".
$cellB7 = $ws.Range("B7")
$cellB7.NumberFormat = "@"
$cellB7.Value = "0"
$cellB7.Style = "Normal"

# D7 (PRCNT_Okanogan_Reaches_EDT_greater_than_RTT for riparian) drops
# from 0.94 to 0.
$ws.Range("D7").Value = 0
